# semana 28 de 2025
# Adds week 28 ("28") as a new column (AE) to the weekly IRA report sheet,
# and fills in a previously-missing value in column AC (week 26) for row 28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell AE1: "28" (stored as text, like the other week-number headers) ---
$ws.Cells.Item(1, 31).Value = "'28"

# --- New values for column AE (week 28), keyed by row number ---
$weekAE = @{
  2  = 52
  4  = 0
  5  = 1
  6  = 59
  7  = 30
  8  = 20
  10 = 1
  12 = 3
  13 = 1
  14 = 2
  15 = 2
  16 = 1
  17 = 4
  18 = 2
  22 = 4
  23 = 6
  24 = 1
  25 = 63
  26 = 2
  28 = 275
  29 = 0
  30 = 22
  31 = 3
  32 = 5
  34 = 0
  35 = 36
  36 = 2
  37 = 10
  38 = 186
  39 = 3
  40 = 36
  41 = 61
  42 = 26
  43 = 201
  44 = 85
  45 = 156
  46 = 3
  47 = 129
  48 = 4
  49 = 0
  50 = 3
  52 = 53
  53 = 0
  54 = 0
  55 = 4
  56 = 13
  57 = 36
}

foreach ($row in $weekAE.Keys) {
  $ws.Cells.Item($row, 31).Value = $weekAE[$row]
}

# --- Fill in the previously-missing week 26 (column AC) value for row 28 ---
$ws.Cells.Item(28, 29).Value = 118
